$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 16,20
$data[0,0] = "ECs"
$data[0,1] = "Vtn"
$data[0,2] = "Itgb3"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 2.767552
$data[0,7] = 8.302655999999999
$data[0,8] = 0.04706493447833917
$data[0,9] = 0.04706493447833917
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 2.759544333333333
$data[0,13] = 8.278633
$data[0,14] = 0.2574067337278401
$data[0,15] = 0.2574067337278401
$data[0,16] = 7.637182438805332
$data[0,17] = 68.73464194924799
$data[0,18] = 0.01211483105718409
$data[0,19] = 0.01211483105718409
$data[1,0] = "ECs"
$data[1,1] = "Vtn"
$data[1,2] = "Itgb3"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 2.767552
$data[1,7] = 8.302655999999999
$data[1,8] = 0.04706493447833917
$data[1,9] = 0.04706493447833917
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 7.245227
$data[1,13] = 21.735681
$data[1,14] = 0.6758254232987829
$data[1,15] = 0.6758254232987829
$data[1,16] = 20.051542474304
$data[1,17] = 180.463882268736
$data[1,18] = 0.03180767926635305
$data[1,19] = 0.03180767926635305
$data[2,0] = "ECs"
$data[2,1] = "Vtn"
$data[2,2] = "Itgb3"
$data[2,3] = "MuSCs"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 2.767552
$data[2,7] = 8.302655999999999
$data[2,8] = 0.04706493447833917
$data[2,9] = 0.04706493447833917
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 0.5200313333333334
$data[2,13] = 1.560094
$data[2,14] = 0.0485078515798926
$data[2,15] = 0.0485078515798926
$data[2,16] = 1.439213756629333
$data[2,17] = 12.952923809664
$data[2,18] = 0.002283018856292647
$data[2,19] = 0.002283018856292646
$data[3,0] = "ECs"
$data[3,1] = "Vtn"
$data[3,2] = "Itgb3"
$data[3,3] = "Resolving-Mac"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 2.767552
$data[3,7] = 8.302655999999999
$data[3,8] = 0.04706493447833917
$data[3,9] = 0.04706493447833917
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 0.1957573333333333
$data[3,13] = 0.587272
$data[3,14] = 0.01825999139348442
$data[3,15] = 0.01825999139348442
$data[3,16] = 0.5417685993813333
$data[3,17] = 4.875917394431999
$data[3,18] = 0.0008594052985093814
$data[3,19] = 0.0008594052985093814
$data[4,0] = "FAPs"
$data[4,1] = "Vtn"
$data[4,2] = "Itgb3"
$data[4,3] = "ECs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 21.05317333333333
$data[4,7] = 63.15952
$data[4,8] = 0.3580298485789791
$data[4,9] = 0.3580298485789791
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 2.759544333333333
$data[4,13] = 8.278633
$data[4,14] = 0.2574067337278401
$data[4,15] = 0.2574067337278401
$data[4,16] = 58.09716517068443
$data[4,17] = 522.8744865361599
$data[4,18] = 0.09215929389978818
$data[4,19] = 0.09215929389978818
$data[5,0] = "FAPs"
$data[5,1] = "Vtn"
$data[5,2] = "Itgb3"
$data[5,3] = "FAPs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 21.05317333333333
$data[5,7] = 63.15952
$data[5,8] = 0.3580298485789791
$data[5,9] = 0.3580298485789791
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 7.245227
$data[5,13] = 21.735681
$data[5,14] = 0.6758254232987829
$data[5,15] = 0.6758254232987829
$data[5,16] = 152.5350198703467
$data[5,17] = 1372.81517883312
$data[5,18] = 0.2419656739694877
$data[5,19] = 0.2419656739694877
$data[6,0] = "FAPs"
$data[6,1] = "Vtn"
$data[6,2] = "Itgb3"
$data[6,3] = "MuSCs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 21.05317333333333
$data[6,7] = 63.15952
$data[6,8] = 0.3580298485789791
$data[6,9] = 0.3580298485789791
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 0.5200313333333334
$data[6,13] = 1.560094
$data[6,14] = 0.0485078515798926
$data[6,15] = 0.0485078515798926
$data[6,16] = 10.94830979943111
$data[6,17] = 98.53478819488001
$data[6,18] = 0.01736725875604054
$data[6,19] = 0.01736725875604054
$data[7,0] = "FAPs"
$data[7,1] = "Vtn"
$data[7,2] = "Itgb3"
$data[7,3] = "Resolving-Mac"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 21.05317333333333
$data[7,7] = 63.15952
$data[7,8] = 0.3580298485789791
$data[7,9] = 0.3580298485789791
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 0.1957573333333333
$data[7,13] = 0.587272
$data[7,14] = 0.01825999139348442
$data[7,15] = 0.01825999139348442
$data[7,16] = 4.121313069937778
$data[7,17] = 37.09181762944
$data[7,18] = 0.00653762195366269
$data[7,19] = 0.00653762195366269
$data[8,0] = "MuSCs"
$data[8,1] = "Vtn"
$data[8,2] = "Itgb3"
$data[8,3] = "ECs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 34.97741266666667
$data[8,7] = 104.932238
$data[8,8] = 0.5948251867999219
$data[8,9] = 0.5948251867999219
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 2.759544333333333
$data[8,13] = 8.278633
$data[8,14] = 0.2574067337278401
$data[8,15] = 0.2574067337278401
$data[8,16] = 96.52172091896156
$data[8,17] = 868.695488270654
$data[8,18] = 0.1531120084732202
$data[8,19] = 0.1531120084732202
$data[9,0] = "MuSCs"
$data[9,1] = "Vtn"
$data[9,2] = "Itgb3"
$data[9,3] = "FAPs"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 34.97741266666667
$data[9,7] = 104.932238
$data[9,8] = 0.5948251867999219
$data[9,9] = 0.5948251867999219
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 7.245227
$data[9,13] = 21.735681
$data[9,14] = 0.6758254232987829
$data[9,15] = 0.6758254232987829
$data[9,16] = 253.4192946426754
$data[9,17] = 2280.773651784078
$data[9,18] = 0.4019979836578348
$data[9,19] = 0.4019979836578348
$data[10,0] = "MuSCs"
$data[10,1] = "Vtn"
$data[10,2] = "Itgb3"
$data[10,3] = "MuSCs"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 34.97741266666667
$data[10,7] = 104.932238
$data[10,8] = 0.5948251867999219
$data[10,9] = 0.5948251867999219
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 0.5200313333333334
$data[10,13] = 1.560094
$data[10,14] = 0.0485078515798926
$data[10,15] = 0.0485078515798926
$data[10,16] = 18.18935054559689
$data[10,17] = 163.704154910372
$data[10,18] = 0.02885369187727251
$data[10,19] = 0.0288536918772725
$data[11,0] = "MuSCs"
$data[11,1] = "Vtn"
$data[11,2] = "Itgb3"
$data[11,3] = "Resolving-Mac"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 34.97741266666667
$data[11,7] = 104.932238
$data[11,8] = 0.5948251867999219
$data[11,9] = 0.5948251867999219
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 0.1957573333333333
$data[11,13] = 0.587272
$data[11,14] = 0.01825999139348442
$data[11,15] = 0.01825999139348442
$data[11,16] = 6.847085030526224
$data[11,17] = 61.62376527473601
$data[11,18] = 0.01086150279159434
$data[11,19] = 0.01086150279159434
$data[12,0] = "Resolving-Mac"
$data[12,1] = "Vtn"
$data[12,2] = "Itgb3"
$data[12,3] = "ECs"
$data[12,4] = 1
$data[12,5] = 0.3333333333333333
$data[12,6] = 0.004706
$data[12,7] = 0.014118
$data[12,8] = 0.00008003014275976175
$data[12,9] = 0.00008003014275976175
$data[12,10] = 3
$data[12,11] = 1
$data[12,12] = 2.759544333333333
$data[12,13] = 8.278633
$data[12,14] = 0.2574067337278401
$data[12,15] = 0.2574067337278401
$data[12,16] = 0.01298641563266667
$data[12,17] = 0.116877740694
$data[12,18] = 0.00002060029764756302
$data[12,19] = 0.00002060029764756302
$data[13,0] = "Resolving-Mac"
$data[13,1] = "Vtn"
$data[13,2] = "Itgb3"
$data[13,3] = "FAPs"
$data[13,4] = 1
$data[13,5] = 0.3333333333333333
$data[13,6] = 0.004706
$data[13,7] = 0.014118
$data[13,8] = 0.00008003014275976175
$data[13,9] = 0.00008003014275976175
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 7.245227
$data[13,13] = 21.735681
$data[13,14] = 0.6758254232987829
$data[13,15] = 0.6758254232987829
$data[13,16] = 0.034096038262
$data[13,17] = 0.306864344358
$data[13,18] = 0.00005408640510727802
$data[13,19] = 0.00005408640510727802
$data[14,0] = "Resolving-Mac"
$data[14,1] = "Vtn"
$data[14,2] = "Itgb3"
$data[14,3] = "MuSCs"
$data[14,4] = 1
$data[14,5] = 0.3333333333333333
$data[14,6] = 0.004706
$data[14,7] = 0.014118
$data[14,8] = 0.00008003014275976175
$data[14,9] = 0.00008003014275976175
$data[14,10] = 3
$data[14,11] = 1
$data[14,12] = 0.5200313333333334
$data[14,13] = 1.560094
$data[14,14] = 0.0485078515798926
$data[14,15] = 0.0485078515798926
$data[14,16] = 0.002447267454666667
$data[14,17] = 0.022025407092
$data[14,18] = 0.00000388209028690814
$data[14,19] = 0.000003882090286908139
$data[15,0] = "Resolving-Mac"
$data[15,1] = "Vtn"
$data[15,2] = "Itgb3"
$data[15,3] = "Resolving-Mac"
$data[15,4] = 1
$data[15,5] = 0.3333333333333333
$data[15,6] = 0.004706
$data[15,7] = 0.014118
$data[15,8] = 0.00008003014275976175
$data[15,9] = 0.00008003014275976175
$data[15,10] = 3
$data[15,11] = 1
$data[15,12] = 0.1957573333333333
$data[15,13] = 0.587272
$data[15,14] = 0.01825999139348442
$data[15,15] = 0.01825999139348442
$data[15,16] = 0.0009212340106666667
$data[15,17] = 0.008291106096
$data[15,18] = 0.000001461349718012579
$data[15,19] = 0.000001461349718012579

$ws.Range("A2:T17").Value = $data
Write-Host "done"
